# Zin.xlsx -- add "ZIN ABS COPY" (L) / "ZIN PHA COPY" (M) columns
# as static snapshots of the existing ZIN ABS (F) / ZIN PHA (G) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers (new shared strings "ZIN ABS COPY" / "ZIN PHA COPY")
$ws.Range("L1").Value = "ZIN ABS COPY"
$ws.Range("M1").Value = "ZIN PHA COPY"

# Static values copied from columns F (ZIN ABS) / G (ZIN PHA) for rows 2-101
$zinAbsCopy = @{}
$zinPhaCopy = @{}
$zinAbsCopy[2] = 1077712.6598583618
$zinPhaCopy[2] = -81.29400291231794
$zinAbsCopy[3] = 1077712.6598583618
$zinPhaCopy[3] = -81.29400291231794
$zinAbsCopy[4] = 1077712.6598583618
$zinPhaCopy[4] = -81.29400291231794
$zinAbsCopy[5] = 1415790.0821762078
$zinPhaCopy[5] = -74.61061750252712
$zinAbsCopy[6] = 1077712.6598583618
$zinPhaCopy[6] = -81.29400291231794
$zinAbsCopy[7] = 956735.2391795018
$zinPhaCopy[7] = -77.17784653260644
$zinAbsCopy[8] = 1073100.9823659612
$zinPhaCopy[8] = -76.65096191642061
$zinAbsCopy[9] = 1415790.0821762078
$zinPhaCopy[9] = -74.61061750252712
$zinAbsCopy[10] = 1073100.9823659612
$zinPhaCopy[10] = -76.65096191642061
$zinAbsCopy[11] = 1073100.9823659612
$zinPhaCopy[11] = -76.65096191642061
$zinAbsCopy[12] = 1073100.9823659612
$zinPhaCopy[12] = -76.65096191642061
$zinAbsCopy[13] = 956735.2391795018
$zinPhaCopy[13] = -77.17784653260644
$zinAbsCopy[14] = 956735.2391795018
$zinPhaCopy[14] = -77.17784653260644
$zinAbsCopy[15] = 956735.2391795018
$zinPhaCopy[15] = -77.17784653260644
$zinAbsCopy[16] = 862995.8160608333
$zinPhaCopy[16] = -77.50365918319714
$zinAbsCopy[17] = 862995.8160608333
$zinPhaCopy[17] = -77.50365918319714
$zinAbsCopy[18] = 862995.8160608333
$zinPhaCopy[18] = -77.50365918319714
$zinAbsCopy[19] = 785918.4832633248
$zinPhaCopy[19] = -77.68202252877072
$zinAbsCopy[20] = 785918.4832633248
$zinPhaCopy[20] = -77.68202252877072
$zinAbsCopy[21] = 721453.2355761374
$zinPhaCopy[21] = -77.74912143961949
$zinAbsCopy[22] = 720242.1188419461
$zinPhaCopy[22] = -74.67016985927062
$zinAbsCopy[23] = 720242.1188419461
$zinPhaCopy[23] = -74.67016985927062
$zinAbsCopy[24] = 619939.3997286303
$zinPhaCopy[24] = -74.99477110740716
$zinAbsCopy[25] = 619939.3997286303
$zinPhaCopy[25] = -74.99477110740716
$zinAbsCopy[26] = 619939.3997286303
$zinPhaCopy[26] = -74.99477110740716
$zinAbsCopy[27] = 579555.4281248677
$zinPhaCopy[27] = -75.02768340099931
$zinAbsCopy[28] = 579555.4281248677
$zinPhaCopy[28] = -75.02768340099931
$zinAbsCopy[29] = 619939.3997286303
$zinPhaCopy[29] = -74.99477110740716
$zinAbsCopy[30] = 579555.4281248677
$zinPhaCopy[30] = -75.02768340099931
$zinAbsCopy[31] = 579039.0917740292
$zinPhaCopy[31] = -72.58110416104186
$zinAbsCopy[32] = 544052.5929942008
$zinPhaCopy[32] = -72.69708653338446
$zinAbsCopy[33] = 544052.5929942008
$zinPhaCopy[33] = -72.69708653338446
$zinAbsCopy[34] = 513040.37608930975
$zinPhaCopy[34] = -72.74285714496085
$zinAbsCopy[35] = 512596.2572960603
$zinPhaCopy[35] = -70.60198079420499
$zinAbsCopy[36] = 512596.2572960603
$zinPhaCopy[36] = -70.60198079420499
$zinAbsCopy[37] = 460728.4048971434
$zinPhaCopy[37] = -70.74346551034604
$zinAbsCopy[38] = 460728.4048971434
$zinPhaCopy[38] = -70.74346551034604
$zinAbsCopy[39] = 418564.04312584404
$zinPhaCopy[39] = -68.94785911776346
$zinAbsCopy[40] = 400367.51010283077
$zinPhaCopy[40] = -68.92838864801526
$zinAbsCopy[41] = 400367.51010283077
$zinPhaCopy[41] = -68.92838864801526
$zinAbsCopy[42] = 383890.33045937447
$zinPhaCopy[42] = -67.3030340776751
$zinAbsCopy[43] = 368786.3754077163
$zinPhaCopy[43] = -65.78382391834847
$zinAbsCopy[44] = 354897.9248741901
$zinPhaCopy[44] = -65.77313506117332
$zinAbsCopy[45] = 354897.9248741901
$zinPhaCopy[45] = -65.77313506117332
$zinAbsCopy[46] = 342168.9389242908
$zinPhaCopy[46] = -64.36296522672225
$zinAbsCopy[47] = 330403.3418057293
$zinPhaCopy[47] = -63.03084620973202
$zinAbsCopy[48] = 319458.20847064734
$zinPhaCopy[48] = -63.013260610872294
$zinAbsCopy[49] = 319505.5571656429
$zinPhaCopy[49] = -61.76827554763051
$zinAbsCopy[50] = 309391.0032417138
$zinPhaCopy[50] = -60.567920882301635
$zinAbsCopy[51] = 299985.36966070416
$zinPhaCopy[51] = -60.56585605426695
$zinAbsCopy[52] = 291232.5665223838
$zinPhaCopy[52] = -59.42566850394155
$zinAbsCopy[53] = 291220.65009764605
$zinPhaCopy[53] = -58.32930101164254
$zinAbsCopy[54] = 283039.68182531523
$zinPhaCopy[54] = -57.280708940100794
$zinAbsCopy[55] = 275389.62584224035
$zinPhaCopy[55] = -56.273447050112274
$zinAbsCopy[56] = 275286.5689678362
$zinPhaCopy[56] = -55.27231430027631
$zinAbsCopy[57] = 268131.154573172
$zinPhaCopy[57] = -54.33986439300911
$zinAbsCopy[58] = 261420.02709892963
$zinPhaCopy[58] = -53.439479137301106
$zinAbsCopy[59] = 255115.59501414662
$zinPhaCopy[59] = -52.56847581231762
$zinAbsCopy[60] = 249184.26046019953
$zinPhaCopy[60] = -51.72447605776703
$zinAbsCopy[61] = 243503.4632357073
$zinPhaCopy[61] = -50.07902167440885
$zinAbsCopy[62] = 238248.6756011388
$zinPhaCopy[62] = -49.309952100789424
$zinAbsCopy[63] = 233186.00112305733
$zinPhaCopy[63] = -47.79636925524512
$zinAbsCopy[64] = 228510.85725144425
$zinPhaCopy[64] = -47.089415776983
$zinAbsCopy[65] = 228243.93985114034
$zinPhaCopy[65] = -45.63576193116603
$zinAbsCopy[66] = 223861.9615291244
$zinPhaCopy[66] = -44.98878020580985
$zinAbsCopy[67] = 219573.4835043044
$zinPhaCopy[67] = -43.683341329922484
$zinAbsCopy[68] = 215516.03914336726
$zinPhaCopy[68] = -42.43675573917941
$zinAbsCopy[69] = 215353.96337797111
$zinPhaCopy[69] = -41.801365050088556
$zinAbsCopy[70] = 211516.689680722
$zinPhaCopy[70] = -40.634854108757516
$zinAbsCopy[71] = 211133.74729009726
$zinPhaCopy[71] = -39.43925876696137
$zinAbsCopy[72] = 207518.33146074356
$zinPhaCopy[72] = -38.36910349978284
$zinAbsCopy[73] = 201199.02729625415
$zinPhaCopy[73] = -37.41519102994378
$zinAbsCopy[74] = 203684.77708271
$zinPhaCopy[74] = -36.267399402807214
$zinAbsCopy[75] = 200474.93115692554
$zinPhaCopy[75] = -35.32069127541109
$zinAbsCopy[76] = 197438.25906580538
$zinPhaCopy[76] = -34.40874867018301
$zinAbsCopy[77] = 194563.6896709653
$zinPhaCopy[77] = -33.529269283091814
$zinAbsCopy[78] = 194178.7419580363
$zinPhaCopy[78] = -32.59665509566369
$zinAbsCopy[79] = 191483.99902355304
$zinPhaCopy[79] = -31.781163600321378
$zinAbsCopy[80] = 191091.238433896
$zinPhaCopy[80] = -30.906062679763846
$zinAbsCopy[81] = 190666.26599825255
$zinPhaCopy[81] = -30.054535301770983
$zinAbsCopy[82] = 187969.42436043106
$zinPhaCopy[82] = -28.923836219577574
$zinAbsCopy[83] = 187538.31041102414
$zinPhaCopy[83] = -28.13517866873207
$zinAbsCopy[84] = 185228.08924077768
$zinPhaCopy[84] = -27.46970558009751
$zinAbsCopy[85] = 184594.3927169932
$zinPhaCopy[85] = -26.36545025362869
$zinAbsCopy[86] = 182454.30997578308
$zinPhaCopy[86] = -25.755438354686625
$zinAbsCopy[87] = 182045.32152795314
$zinPhaCopy[87] = -25.067630514915233
$zinAbsCopy[88] = 181622.29286118617
$zinPhaCopy[88] = -24.398545060277737
$zinAbsCopy[89] = 181187.0825153705
$zinPhaCopy[89] = -23.74777319996105
$zinAbsCopy[90] = 182242.5782938723
$zinPhaCopy[90] = -23.001635092538308
$zinAbsCopy[91] = 180056.8327893026
$zinPhaCopy[91] = -22.19823174683693
$zinAbsCopy[92] = 179591.84452269485
$zinPhaCopy[92] = -21.608265815022435
$zinAbsCopy[93] = 179121.54771915384
$zinPhaCopy[93] = -21.03472445902592
$zinAbsCopy[94] = 181268.47422416843
$zinPhaCopy[94] = -20.227632693364047
$zinAbsCopy[95] = 180709.42584265923
$zinPhaCopy[95] = -19.684373438273923
$zinAbsCopy[96] = 179872.16037297508
$zinPhaCopy[96] = -18.898669081794246
$zinAbsCopy[97] = 179316.15967160516
$zinPhaCopy[97] = -18.393674702711838
$zinAbsCopy[98] = 178762.83206885803
$zinPhaCopy[98] = -17.903219325221706
$zinAbsCopy[99] = 178212.88794370918
$zinPhaCopy[99] = -17.4268773221811
$zinAbsCopy[100] = 177666.96514221444
$zinPhaCopy[100] = -16.96423213088691
$zinAbsCopy[101] = 178165.28777603325
$zinPhaCopy[101] = -16.37903807463735

for ($r = 2; $r -le 101; $r++) {
    $ws.Cells.Item($r, 12).Value = $zinAbsCopy[$r]
    $ws.Cells.Item($r, 13).Value = $zinPhaCopy[$r]
}

# Match the author's final selection/scroll state: cell L2 selected, no frozen scroll offset
$null = $ws.Range("L2").Select()
